# [Kadastro App] Yeni kayit eklendi: 2942
# Append the new record to both the "Kayitlar" (master) sheet and the
# related "Erdemli" district sheet.

$wb = $excel.ActiveWorkbook

$newRow = @("2942", "2025-09-08", "Erdemli", "1", "MAHKEME KARARI", "EMİNE ALANLI KIRCILI (K.Mühendisi), SERDAR ARSLAN (Tekniker)")

$sheetNames = @("Kayitlar", "Erdemli")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Find first empty row below the existing data in column A.
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    $targetRow = $lastRow + 1

    for ($col = 1; $col -le $newRow.Count; $col++) {
        $cell = $ws.Cells.Item($targetRow, $col)
        # Force text storage (so numeric-looking values like "2942" or
        # "2025-09-08" are kept as text, matching the rest of the sheet)
        # without leaving the cell tagged with a non-default number format.
        $cell.NumberFormat = "@"
        $cell.Value = $newRow[$col - 1]
        $cell.Style = "Normal"
    }
}
